# Apply the latest cryptos-list snapshot (price + 1h-volume refresh,
# plus a handful of coins whose rank/row order changed) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.489.38"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.796.66"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'316.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5421"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.3790"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'0.07506"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.111"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'41.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'20.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.97%  "
$ws.Range("D14").Value = "'6.165"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "'7.314"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "1.793.59"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "'89.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "'0.00001067"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'0.06484"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "'17.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "'5.971"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "28.476.88"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'11.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'2.078"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "'159.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("D27").Value = "'20.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "1.999.80"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "'2.312"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.36%  "
$ws.Range("D30").Value = "'122.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").Value = "'1.109"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").Value = "'5.624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "'0.2278"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "'0.06479"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.38%  "
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.036"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.596"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("D40").Value = "'11.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.96%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.198"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.6197"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.450"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.97%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "'3.683"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "'0.5819"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("D48").Value = "'127.52"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.950"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.197"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("D51").Value = "'0.06873"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.42%  "
